# Update cryptocurrency price list (Price / Volume(1h) columns) with the
# latest snapshot values, preserving the original text-cell formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric but must remain stored as text
$textFormatCells = @(
    'D5',
    'D6',
    'D7',
    'D8',
    'D10',
    'D13',
    'D16',
    'D17',
    'D19',
    'D20',
    'D21',
    'D22',
    'D24',
    'D26',
    'D28',
    'D31',
    'D32',
    'D34',
    'D36',
    'D37',
    'D39',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D48',
    'D49',
    'D50',
    'D51'
)

foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# New cell values taken from the updated cryptos snapshot
$updates = [ordered]@{
    'D2' = '58.703.86'
    'E2' = '  -0.01%  '
    'D3' = '2.569.32'
    'E3' = '  -0.66%  '
    'E4' = '  +0.11%  '
    'D5' = '562.07'
    'E5' = '  +1.96%  '
    'D6' = '142.90'
    'E6' = '  -1.14%  '
    'D7' = '1.00'
    'E7' = '  +0.19%  '
    'D8' = '0.597'
    'E8' = '  +0.53%  '
    'D9' = '2.576.33'
    'E9' = '  -0.43%  '
    'D10' = '6.65'
    'E10' = '  -2.74%  '
    'E11' = '  +2.26%  '
    'E12' = '  +8.47%  '
    'D13' = '0.341'
    'E13' = '  +2.29%  '
    'D14' = '3.023.53'
    'E14' = '  -0.37%  '
    'D15' = '58.807.37'
    'E15' = '  +0.32%  '
    'D16' = '22.02'
    'E16' = '  +6.27%  '
    'D17' = '0.0000136'
    'E17' = '  +3.79%  '
    'D18' = '2.572.33'
    'E18' = '  -0.45%  '
    'D19' = '4.48'
    'E19' = '  +1.03%  '
    'D20' = '334.83'
    'E20' = '  -0.23%  '
    'D21' = '10.14'
    'E21' = '  +1.12%  '
    'D22' = '6.14'
    'E22' = '  +0.60%  '
    'E23' = '  +0.09%  '
    'D24' = '63.67'
    'E24' = '  -4.18%  '
    'E25' = '  +6.28%  '
    'D26' = '0.998'
    'E26' = '  -0.02%  '
    'E27' = '  +1.79%  '
    'D28' = '7.24'
    'E28' = '  +2.39%  '
    'D29' = '0.0₃0775'
    'E29' = '  +4.05%  '
    'E30' = '  +0.07%  '
    'D31' = '1.66'
    'E31' = '  +0.18%  '
    'D32' = '158.24'
    'E32' = '  +2.30%  '
    'E33' = '  +2.39%  '
    'D34' = '18.94'
    'E34' = '  +0.55%  '
    'E35' = '  +2.24%  '
    'D36' = '0.872'
    'E36' = '  +2.16%  '
    'D37' = '0.875'
    'E37' = '  +6.71%  '
    'E38' = '  +1.89%  '
    'D39' = '36.70'
    'E39' = '  -1.39%  '
    'E40' = '  +2.96%  '
    'D41' = '290.04'
    'E41' = '  +4.06%  '
    'D42' = '3.63'
    'E42' = '  +1.05%  '
    'D43' = '0.999'
    'E43' = '  +0.17%  '
    'D44' = '0.0969'
    'E44' = '  +2.15%  '
    'D45' = '0.592'
    'E45' = '  -0.36%  '
    'E46' = '  -0.36%  '
    'E47' = '  +0.63%  '
    'D48' = '19.05'
    'E48' = '  +2.37%  '
    'D49' = '123.72'
    'E49' = '  +10.23%  '
    'D50' = '0.0230'
    'E50' = '  +1.88%  '
    'D51' = '18.45'
    'E51' = '  +3.60%  '
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

